$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Add new column K (year 2023) to the table ---

# Header row (row 3): copy format from I3 for the new K3 cell, then set its value.
$ws.Range("I3").Copy() | Out-Null
$ws.Range("K3").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("K3").Value = 2023

# The existing J3 (year 2022) header cell swaps to the "non-colored" header style,
# matching the style used by B3/C3.
$ws.Range("C3").Copy() | Out-Null
$ws.Range("J3").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("J3").Value = 2022

# Data rows 4-6: copy format from column I (the preceding data column) into K, then
# fill in the new figures for 2023.
$ws.Range("I4").Copy() | Out-Null
$ws.Range("K4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("K4").Value = 1266.7

$ws.Range("I5").Copy() | Out-Null
$ws.Range("K5").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("K5").Value = 867.9

$ws.Range("I6").Copy() | Out-Null
$ws.Range("K6").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("K6").Value = 1444.8
